$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (new Price text, new Volume(1h) text). Price is omitted (left as $null)
# for rows where the diff only changes the Volume column.
$updates = @(
    @{ Row = 2; Price = "42.819.07"; Volume = "  +0.47%  " }
    @{ Row = 3; Price = "2.284.45"; Volume = "  -0.41%  " }
    @{ Row = 4; Price = $null; Volume = "  +0.10%  " }
    @{ Row = 5; Price = "309.34"; Volume = "  -4.15%  " }
    @{ Row = 6; Price = "102.61"; Volume = "  -1.18%  " }
    @{ Row = 7; Price = $null; Volume = "  -1.34%  " }
    @{ Row = 8; Price = $null; Volume = "  +0.03%  " }
    @{ Row = 9; Price = $null; Volume = "  -1.56%  " }
    @{ Row = 10; Price = $null; Volume = "  -4.15%  " }
    @{ Row = 11; Price = $null; Volume = "  -1.03%  " }
    @{ Row = 12; Price = "8.16"; Volume = "  -2.98%  " }
    @{ Row = 13; Price = $null; Volume = "  +0.43%  " }
    @{ Row = 14; Price = $null; Volume = "  -0.81%  " }
    @{ Row = 15; Price = "15.16"; Volume = "  -0.51%  " }
    @{ Row = 16; Price = "2.631.04"; Volume = "  -0.38%  " }
    @{ Row = 17; Price = "2.285.34"; Volume = "  -0.57%  " }
    @{ Row = 18; Price = "42.444.83"; Volume = "  -0.25%  " }
    @{ Row = 19; Price = "7.26"; Volume = "  -2.18%  " }
    @{ Row = 20; Price = $null; Volume = "  -1.61%  " }
    @{ Row = 21; Price = "13.31"; Volume = "  -2.51%  " }
    @{ Row = 22; Price = "73.09"; Volume = "  -0.33%  " }
    @{ Row = 23; Price = "269.94"; Volume = "  +0.12%  " }
    @{ Row = 24; Price = "3.37"; Volume = "  -6.12%  " }
    @{ Row = 25; Price = $null; Volume = "  -3.40%  " }
    @{ Row = 26; Price = $null; Volume = "  -0.40%  " }
    @{ Row = 27; Price = "10.69"; Volume = "  -2.14%  " }
    @{ Row = 28; Price = "6.93"; Volume = "  +12.13%  " }
    @{ Row = 29; Price = $null; Volume = "  -2.33%  " }
    @{ Row = 30; Price = "22.31"; Volume = "  -1.06%  " }
    @{ Row = 31; Price = $null; Volume = "  -6.87%  " }
    @{ Row = 32; Price = "164.25"; Volume = "  -0.52%  " }
    @{ Row = 33; Price = "0.0845"; Volume = "  -4.26%  " }
    @{ Row = 34; Price = $null; Volume = "  -2.86%  " }
    @{ Row = 35; Price = $null; Volume = "  +0.93%  " }
    @{ Row = 36; Price = "0.111"; Volume = "  -3.28%  " }
    @{ Row = 37; Price = $null; Volume = "  -3.34%  " }
    @{ Row = 38; Price = $null; Volume = "  -3.24%  " }
    @{ Row = 39; Price = $null; Volume = "  +0.39%  " }
    @{ Row = 40; Price = $null; Volume = "  -2.89%  " }
    @{ Row = 41; Price = "112.02"; Volume = "  +20.37%  " }
    @{ Row = 42; Price = $null; Volume = "  +0.91%  " }
    @{ Row = 43; Price = "69.11"; Volume = "  -0.68%  " }
    @{ Row = 44; Price = $null; Volume = "  -0.33%  " }
    @{ Row = 45; Price = $null; Volume = "  -0.77%  " }
    @{ Row = 46; Price = $null; Volume = "  -2.82%  " }
    @{ Row = 47; Price = "1.703.67"; Volume = "  +6.42%  " }
    @{ Row = 48; Price = "109.94"; Volume = "  -2.89%  " }
    @{ Row = 49; Price = "77.02"; Volume = "  -5.90%  " }
    @{ Row = 50; Price = "5.16"; Volume = "  -2.25%  " }
    @{ Row = 51; Price = $null; Volume = "  -3.60%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.Price) {
        $priceCell = $ws.Cells.Item($u.Row, 4)
        # Force text storage so numeric-looking prices (e.g. "309.34")
        # are not auto-coerced into the Number type by Excel, then restore
        # the original (default) cell style so no stray formatting is left behind.
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $u.Price
        $priceCell.Style = "Normal"
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.Volume
}
